$wb = $excel.ActiveWorkbook

# --- 1. Insert a new worksheet "6_1_4" right after "6_1_3" (before "6_1_5") ---
$ws3 = $wb.Worksheets.Item("6_1_3")
$new = $wb.Worksheets.Add($null, $ws3)
$new.Name = "6_1_4"

$new.Range("B1").Value = 0
$new.Range("A2").Value = 0
$new.Range("B2").Value = "intersection"
$new.Range("A3").Value = 1
$new.Range("B3").Value = 51

# Match the bold/bordered/centered style used by the sibling summary sheets (e.g. 6_1_2!B1)
$styleSource = $wb.Worksheets.Item("6_1_2").Range("B1")
$styleSource.Copy()
$new.Range("B1").PasteSpecial(-4122)
$new.Range("A2").PasteSpecial(-4122)
$new.Range("A3").PasteSpecial(-4122)

# --- 2. Update the Janre / Players_median sheet (still named "6_1_5") ---
# Rows 2-5 get cyclically rotated down by one (row2 <- old row5, row3 <- old row2, etc.)
$janre = $wb.Worksheets.Item("6_1_5")

$janre.Range("B2").Value = "Party"
$janre.Range("C2").Value = 11

$janre.Range("B3").Value = "Мафия"
$janre.Range("C3").Value = 8.5

$janre.Range("B4").Value = "Карточная"
$janre.Range("C4").Value = 4.5

$janre.Range("B5").Value = "Семейная"
$janre.Range("C5").Value = 5
